# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 ("R40") becomes the text "1" (leading apostrophe forces a text
# entry, matching the shared-string / t="s" cell the diff expects rather
# than a numeric 1).
$ws.Range("B11").Value = "'1"
